$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new bitcoin-buy entry for 2025-09-03 as row 46, mirroring the
# existing rows where the Date column holds plain text (e.g. "08/31/2025")
# rather than a date serial value.
$row = 46

$dateCell = $ws.Cells.Item($row, 1)
# Force text interpretation so "09/03/2025" is stored as a string instead
# of being auto-converted into a date serial number, then drop back to the
# default "Normal" style so the cell carries no explicit style index.
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/03/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.0004422799999999998
$ws.Cells.Item($row, 3).Value = 113050.5562087366
$ws.Cells.Item($row, 4).Value = 50
